$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "ROW50-FE-LIFTER" (sheet1): append row 95
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$ws1.Range("A95").NumberFormat = $ws1.Range("A94").NumberFormat
$ws1.Range("A95").Value = 45771.80362907407
$ws1.Range("B95").Value = "0x01,0x90"
$ws1.Range("C95").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Range("D95").Value = "0x01,0x3a"
$ws1.Range("E95").Value = "0xe"
$ws1.Range("F95").Value = 400
$ws1.Range("G95").Value = 568631262647114000000000.0
$ws1.Range("H95").Value = 314
$ws1.Range("I95").Value = 14

# ---------------------------------------------------------------------------
# Sheet "ROW50-MID-LIFTER" (sheet2): append row 97
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$ws2.Range("A97").NumberFormat = $ws2.Range("A96").NumberFormat
$ws2.Range("A97").Value = 45771.76709490741
$ws2.Range("B97").Value = "0x01,0x90 "
$ws2.Range("C97").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Range("D97").Value = "0x01,0x42"
$ws2.Range("E97").Value = "0x19"
$ws2.Range("F97").Value = 400
# G97 stores a text representation of the big integer (too large for a double
# to round-trip) - force text storage without leaving the cell's own style
# dirtied (ClearFormats resets the per-cell style index back to default after
# the text number-format coercion has taken effect).
$ws2.Range("G97").NumberFormat = "@"
$ws2.Range("G97").Value = "568631262647113771663628"
$ws2.Range("G97").ClearFormats()
$ws2.Range("H97").Value = 322
$ws2.Range("I97").Value = 25

# ---------------------------------------------------------------------------
# Sheet "ROW11-FE-LIFTER" (sheet3): append row 95
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$ws3.Range("A95").NumberFormat = $ws3.Range("A94").NumberFormat
$ws3.Range("A95").Value = 45771.83461414352
$ws3.Range("B95").Value = "0x01,0x90"
$ws3.Range("C95").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Range("D95").Value = "0x01,0x3a"
$ws3.Range("E95").Value = "0x14"
$ws3.Range("F95").Value = 400
$ws3.Range("G95").Value = 568631262647114000000000.0
$ws3.Range("H95").Value = 314
$ws3.Range("I95").Value = 20

# ---------------------------------------------------------------------------
# Sheet "ROW11-MID-LIFTER" (sheet4): append row 95
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$ws4.Range("A95").NumberFormat = $ws4.Range("A94").NumberFormat
$ws4.Range("A95").Value = 45771.95616605324
$ws4.Range("B95").Value = "0x01,0x90"
$ws4.Range("C95").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Range("D95").Value = "0x01,0x42"
$ws4.Range("E95").Value = "0x19"
$ws4.Range("F95").Value = 400
$ws4.Range("G95").Value = 568631262647114000000000.0
$ws4.Range("H95").Value = 322
$ws4.Range("I95").Value = 25
